$d = $word.ActiveDocument

# --- Locate the target paragraph: the bullet ending with
# "...government and DHS, they will create new standardized education to fill these jobs." ---
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text
    if ($ptxt -like "*they will create new standardized education to fill these jobs.*") {
        $targetParaIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetParaIndex)
$pr = $p.Range
$pr.End = $pr.End - 1
$lastRunText = "government and DHS, they will create new standardized education to fill these jobs."
$fullText = $pr.Text
$idx = $fullText.LastIndexOf($lastRunText)
$runStart = $pr.Start + $idx
$runEnd = $pr.End

# Remove the old final run's text, then re-insert it (extended with " Further")
# as a fresh run so it lands as a single run, same as the original.
$runRange = $d.Range($runStart, $runEnd)
$runRange.Delete()

$p = $d.Paragraphs.Item($targetParaIndex)
$pr = $p.Range
$pr.End = $pr.End - 1
$insertPoint = $d.Range($pr.End, $pr.End)
$insertPoint.InsertAfter($lastRunText + " Further")

# --- Append the rest of the new sentences to the same paragraph, each
# line as its own run with single-space runs in between (matching the
# existing run-per-line convention used throughout this document). ---
$moreRuns = @(
    " ",
    "study on good solutions for this problem can be beneficial, a new Cybersecurity certification",
    " ",
    "offered by the federal govenment or making Cybersecurity courses ubiquitous at every education",
    " ",
    "institution at the community college and upperclass high school level; Funnelling students",
    " ",
    "into roles that help the government at the local levels."
)

foreach ($chunk in $moreRuns) {
    $p = $d.Paragraphs.Item($targetParaIndex)
    $pr = $p.Range
    $pr.End = $pr.End - 1
    $ip = $d.Range($pr.End, $pr.End)
    $ip.InsertAfter($chunk)
}

# --- Insert a brand-new bullet paragraph right after it, with the same
# numbering (numId 1000), and fill it in with the same run-per-line style. ---
$p = $d.Paragraphs.Item($targetParaIndex)
$pr = $p.Range
$pr.End = $pr.End - 1
$ip = $d.Range($pr.End, $pr.End)
$ip.InsertParagraphAfter()

$newParaIndex = $targetParaIndex + 1

$newParaRuns = @(
    "No matter the outcome of this bill I can say that S. 2520 will negatively effect black-hat",
    " ",
    "hackers by tighting up the security on low hanging fruit targets they would otherwise hack",
    " ",
    "or not get caught hacking. Penetration testers will now have to be careful when running",
    " ",
    "intrusive enumeration scanning techniques on I.P. address space ranges that might contain",
    " ",
    "hosts to government entities."
)

foreach ($chunk in $newParaRuns) {
    $p2 = $d.Paragraphs.Item($newParaIndex)
    $pr2 = $p2.Range
    $pr2.End = $pr2.End - 1
    $ip2 = $d.Range($pr2.End, $pr2.End)
    $ip2.InsertAfter($chunk)
}
